$d = $word.ActiveDocument
$t = $d.Tables.Item(1)

# Each data row in the table is followed by 3 blank rows; rows 1,5,9,13,17
# hold the five "NN÷N=" expressions that changed in this revision.

$changes = @(
    @{ Row = 1;  Col = 1; Old = "13÷7="; New = "36÷7=" },
    @{ Row = 1;  Col = 2; Old = "98÷7="; New = "74÷2=" },
    @{ Row = 1;  Col = 3; Old = "55÷8="; New = "19÷8=" },
    @{ Row = 1;  Col = 4; Old = "36÷3="; New = "58÷6=" },
    @{ Row = 1;  Col = 5; Old = "43÷3="; New = "91÷2=" },

    @{ Row = 5;  Col = 1; Old = "10÷5="; New = "11÷9=" },
    @{ Row = 5;  Col = 2; Old = "46÷8="; New = "29÷5=" },
    @{ Row = 5;  Col = 3; Old = "79÷2="; New = "80÷4=" },
    @{ Row = 5;  Col = 4; Old = "88÷5="; New = "16÷4=" },
    @{ Row = 5;  Col = 5; Old = "73÷3="; New = "11÷3=" },

    @{ Row = 9;  Col = 1; Old = "13÷8="; New = "86÷5=" },
    @{ Row = 9;  Col = 2; Old = "96÷3="; New = "49÷4=" },
    @{ Row = 9;  Col = 3; Old = "83÷8="; New = "93÷5=" },
    @{ Row = 9;  Col = 4; Old = "29÷5="; New = "43÷3=" },
    @{ Row = 9;  Col = 5; Old = "27÷4="; New = "47÷4=" },

    @{ Row = 13; Col = 1; Old = "45÷7="; New = "77÷6=" },
    @{ Row = 13; Col = 2; Old = "33÷6="; New = "71÷4=" },
    @{ Row = 13; Col = 3; Old = "46÷8="; New = "96÷3=" },
    @{ Row = 13; Col = 4; Old = "19÷4="; New = "92÷2=" },
    @{ Row = 13; Col = 5; Old = "30÷3="; New = "22÷9=" },

    @{ Row = 17; Col = 1; Old = "42÷9="; New = "52÷8=" },
    @{ Row = 17; Col = 2; Old = "63÷3="; New = "71÷8=" },
    @{ Row = 17; Col = 3; Old = "57÷6="; New = "87÷2=" },
    @{ Row = 17; Col = 4; Old = "44÷6="; New = "45÷5=" },
    @{ Row = 17; Col = 5; Old = "43÷9="; New = "33÷8=" }
)

# Several of the expressions (e.g. "46÷8=", "29÷5=") occur more than once in
# the table, and some values created by one change become the "old" value
# looked up by a later change. A document-wide Find/Replace (even one scoped
# to a cell Range) is unsafe here, so each target cell's Range.Text is set
# directly instead, which only touches that single cell and leaves the run
# formatting (rFonts/sz) untouched.
foreach ($change in $changes) {
    $cell = $t.Cell($change.Row, $change.Col)
    $range = $cell.Range
    # Cell.Range.Text includes the trailing end-of-cell mark (CR + BEL), so
    # compare only the leading part that corresponds to the visible text.
    $current = $range.Text.Substring(0, $change.Old.Length)
    if ($current -ne $change.Old) {
        throw "Unexpected text in cell (" + $change.Row + "," + $change.Col + "): [" + $current + "], expected [" + $change.Old + "]"
    }
    $range.Text = $change.New
}
